# OpenTBS 1.7.0-beta, can adjust pictures
#
# This script reproduces, through the Word object model, the two textual
# edits made to word/document.xml:
#
#   1. The old "_GoBack" bookmark that used to split the sentence
#      "... tags are automatically mer|ged in headers and footers."
#      is gone; the two runs around it collapse back into a single run.
#
#   2. The demo tag
#         [b.number;ope=changepic;from=pic_[val].png;default=current]
#      gets ";adjust" typed in just before the closing "]", and Word's
#      "last edit" bookmark ("_GoBack") now marks that new insertion
#      point, i.e. right after ";adjust" and right before "]".
#
# (The chart axis IDs and the customXml GUID that also changed in the
# original commit are internal artifacts Word re-mints on save; they
# are not reachable through the Word object model, so they are left
# untouched here.)

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Part 1: "... automatically merged in headers and footers."
# -----------------------------------------------------------------
# Re-typing the (unchanged) sentence over itself makes the editor
# normalize/merge the two runs that used to be separated by the old
# "_GoBack" bookmark, producing a single run with the full sentence.
$rng1 = $d.Content
$rng1.Find.Execute("automatically merged in headers and footers.", $true, $false, $false, $false, $false, $true, 1, $false, "automatically merged in headers and footers.", 2)

# -----------------------------------------------------------------
# Part 2: "...default=current]" -> "...default=current;adjust]"
# -----------------------------------------------------------------
# Locate the point right after "default=current" (and right before the
# closing "]").
$rng2 = $d.Content
$rng2.Find.Execute("default=current", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$editPoint = $rng2.End

# Type the new text in at that point.
$insertRng = $d.Range($editPoint, $editPoint)
$insertRng.InsertAfter(";adjust")
$afterAdjust = $insertRng.End

# Temporarily mark the original edit point so the run that was typed
# into does not get silently re-merged with the text that follows
# ("]"); this reproduces Word leaving ";adjust" as its own run.
$splitRng = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("TMP_SPLIT_MARK", $splitRng)

# Move Word's "last edit location" bookmark ("_GoBack") to right after
# the newly typed ";adjust", i.e. right before the closing "]".
$goBackRng = $d.Range($afterAdjust, $afterAdjust)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# Remove the temporary helper bookmark; the run split it created stays.
$tmpBookmark = $d.Bookmarks.Item("TMP_SPLIT_MARK")
$tmpBookmark.Delete()
